$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as produced by the
# "Updated cryptos list" GitHub Actions workflow run.

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.963.69"
Set-TextValue "E2" "  +4.09%  "
Set-TextValue "D3" "2.244.56"
Set-TextValue "E3" "  +3.14%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "245.24"
Set-TextValue "E5" "  +3.09%  "
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "D7" "76.02"
Set-TextValue "E7" "  +8.11%  "
Set-TextValue "E8" "  -0.15%  "
Set-TextValue "E9" "  +6.18%  "
Set-TextValue "D10" "41.09"
Set-TextValue "E10" "  +1.78%  "
Set-TextValue "D11" "0.0934"
Set-TextValue "E11" "  +0.68%  "
Set-TextValue "D12" "6.96"
Set-TextValue "E12" "  +2.48%  "
Set-TextValue "E13" "  +0.78%  "
Set-TextValue "D14" "2.561.44"
Set-TextValue "E14" "  +2.51%  "
Set-TextValue "D15" "14.62"
Set-TextValue "E15" "  +4.48%  "
Set-TextValue "D16" "2.253.81"
Set-TextValue "E16" "  +3.47%  "
Set-TextValue "D17" "0.809"
Set-TextValue "E17" "  +0.77%  "
Set-TextValue "D18" "42.866.96"
Set-TextValue "E18" "  +4.32%  "
Set-TextValue "D19" "0.0000105"
Set-TextValue "E19" "  +3.42%  "
Set-TextValue "D20" "71.21"
Set-TextValue "E20" "  +0.94%  "
Set-TextValue "E21" "  +0.88%  "
Set-TextValue "D22" "10.08"
Set-TextValue "E22" "  +2.04%  "
Set-TextValue "D23" "2.21"
Set-TextValue "E23" "  +13.69%  "
Set-TextValue "D24" "231.28"
Set-TextValue "E24" "  +2.17%  "
Set-TextValue "E26" "  -0.15%  "
Set-TextValue "D27" "3.44"
Set-TextValue "E27" "  -3.08%  "
Set-TextValue "E28" "  +2.21%  "
Set-TextValue "D29" "38.53"
Set-TextValue "E29" "  +23.07%  "
Set-TextValue "B30" "Monero"
Set-TextValue "C30" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D30" "173.05"
Set-TextValue "E30" "  +3.19%  "
Set-TextValue "B31" "Toncoin"
Set-TextValue "C31" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D31" "2.15"
Set-TextValue "E31" "  -1.63%  "
Set-TextValue "E32" "  +1.69%  "
Set-TextValue "D33" "0.0798"
Set-TextValue "E33" "  +3.66%  "
Set-TextValue "D34" "5.34"
Set-TextValue "E34" "  +3.21%  "
Set-TextValue "D35" "0.122"
Set-TextValue "E35" "  +1.22%  "
Set-TextValue "E36" "  +6.66%  "
Set-TextValue "E37" "  +4.80%  "
Set-TextValue "D38" "0.0330"
Set-TextValue "E38" "  +15.53%  "
Set-TextValue "D39" "12.83"
Set-TextValue "E39" "  +7.71%  "
Set-TextValue "D40" "2.14"
Set-TextValue "E40" "  +2.78%  "
Set-TextValue "E41" "  +2.55%  "
Set-TextValue "E42" "  +6.71%  "
Set-TextValue "D43" "107.48"
Set-TextValue "E43" "  +9.17%  "
Set-TextValue "D44" "59.81"
Set-TextValue "E44" "  -0.45%  "
Set-TextValue "D45" "8.70"
Set-TextValue "E45" "  +4.86%  "
Set-TextValue "D46" "0.0990"
Set-TextValue "E46" "  +1.25%  "
Set-TextValue "D47" "0.458"
Set-TextValue "E47" "  +20.47%  "
Set-TextValue "E48" "  +5.57%  "
Set-TextValue "E49" "  +1.55%  "
Set-TextValue "E50" "  +1.45%  "
Set-TextValue "D51" "2.449.97"
Set-TextValue "E51" "  +3.06%  "
